$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.538636922836304
$ws.Range("B1").Value = 1.767015814781189
$ws.Range("C1").Value = 2.218263864517212
$ws.Range("D1").Value = 4.607551097869873
$ws.Range("E1").Value = 2.194780588150024
